# Generate Report for Handback
# Updates handback timestamps and status on the Overview, zh-cn and de-de
# sheets to reflect a newer handback run for the
# "4e703979-fdb7-4215-887f-3e0a3a6d7f4f.md" / "a337b7ef-2a17-432d-84b8-3e323557cf5e.md"
# files (they were generated/handed back together, hence sharing values).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-28 16:18:04"
$wsOverview.Range("G4").Value = "2016-08-28 16:18:04"

# --- zh-cn sheet ---
# Column E = "Status", H = "Correspond Handoff Datetime", K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-28 16:17:57"
$wsZhCn.Range("H4").Value = "2016-08-28 16:17:57"
$wsZhCn.Range("K3").Value = "2016-08-28 16:18:25"
$wsZhCn.Range("K4").Value = "2016-08-28 16:18:25"

# --- de-de sheet ---
# Column E = "Status", H = "Correspond Handoff Datetime", K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-28 16:18:04"
$wsDeDe.Range("H4").Value = "2016-08-28 16:18:04"
$wsDeDe.Range("K3").Value = "2016-08-28 16:18:32"
$wsDeDe.Range("K4").Value = "2016-08-28 16:18:32"
